$wb = $excel.ActiveWorkbook

# --- Rename Sheet4 -> MeansForExport, add a new "Notes" sheet after it ---
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Name = "MeansForExport"

$notes = $wb.Worksheets.Add($null, $ws4)
$notes.Name = "Notes"
$notes.Range("A1").Value = "GFP-Rab17 images are actually of GFP-Rab24"
$notes.Range("A2").Value = "GFP-Rab24 images are actually of GFP-Rab17"
$notes.Range("A3").Select()

# --- Swap the GFP-Rab17 (col T) and GFP-Rab24 (col AA) values for rows 2-11 ---
for ($r = 2; $r -le 11; $r++) {
    $tCell = $ws4.Cells.Item($r, 20)   # column T
    $aaCell = $ws4.Cells.Item($r, 27)  # column AA
    $tVal = $tCell.Value2
    $aaVal = $aaCell.Value2

    if ($tVal -eq "" -or $tVal -eq $null) {
        $aaCell.ClearContents()
    } else {
        $aaCell.Value = $tVal
    }

    if ($aaVal -eq "" -or $aaVal -eq $null) {
        $tCell.ClearContents()
    } else {
        $tCell.Value = $aaVal
    }
}

# --- Fix up Sheet3's stale selection scroll (diff drops topLeftCell) ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("B385:B394").Select()

# --- Re-activate MeansForExport as the active sheet/tab with the new selection ---
$ws4.Activate()
$ws4.Range("T2").Select()
